$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.827.33"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "2.089.71"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("E12").Value = "  +2.58%  "

$ws.Range("D13").Value = "2.398.45"
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.780"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("D17").Value = "2.090.76"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "37.837.34"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.35%  "

$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0235"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.61%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.454.45"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("E47").Value = "  -4.37%  "

$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "

$ws.Range("D51").Value = "2.282.29"
$ws.Range("E51").Value = "  +0.01%  "

